# Commit: "Some new cards. VP point calculations"
#
# 1) Rename Sheet1 -> Deck, Sheet2 -> VPs, delete Sheet3.
# 2) Populate the VPs sheet with a Resource / VP Rate lookup table.
# 3) Insert 5 new card rows into the Deck sheet (rows 16-20), pushing the
#    old rows 16-22 down to 21-27.
# 4) Fill in the data for the 5 new cards.
# 5) Add the VP-calculation formula to column G for every card row.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets & drop Sheet3 -----------------------------------
$wb.Worksheets.Item("Sheet1").Name = "Deck"
$wb.Worksheets.Item("Sheet2").Name = "VPs"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet3").Delete()
$excel.DisplayAlerts = $true

$deck = $wb.Worksheets.Item("Deck")
$vps  = $wb.Worksheets.Item("VPs")

# --- 2. Populate the VPs sheet -----------------------------------------
$vps.Range("A1").Value = "Resource"
$vps.Range("B1").Value = "VP Rate"
$vps.Range("A1:B1").Font.Bold = $true

$vps.Range("A2").Value = "Wood"
$vps.Range("B2").Value = 1.3
$vps.Range("A3").Value = "Steel"
$vps.Range("B3").Value = 1.4
$vps.Range("A4").Value = "Stone"
$vps.Range("B4").Value = 1.5
$vps.Range("A5").Value = "Gold"
$vps.Range("B5").Value = 1.6
$vps.Range("B6").Value = 0

$vps.Columns.Item(2).ColumnWidth = 9

# --- 3. Insert 5 new rows at 16-20 in the Deck sheet --------------------
$deck.Range("A16:A20").EntireRow.Insert()
# The inserted rows inherit the old row 16's column-I style (it used to
# have a description cell there); the new cards don't use columns I/J,
# so drop that leftover formatting.
$deck.Range("I16:J20").Clear()

# --- 4. Fill in the new cards --------------------------------------------
# Row 16: Throne
$deck.Range("A16").Value = "Throne"
$deck.Range("B16").Value = 1
$deck.Range("C16").Value = 2
$deck.Range("D16").Value = "Stone"
$deck.Range("E16").Value = 2
$deck.Range("F16").Value = "Gold"
$deck.Range("H16").Value = "stone-throne"

# Row 17: Stone Tablet
$deck.Range("A17").Value = "Stone Tablet"
$deck.Range("B17").Value = 1
$deck.Range("C17").Value = 3
$deck.Range("D17").Value = "Stone"
$deck.Range("H17").Value = "stone-tablet"

# Row 18: Anvil
$deck.Range("A18").Value = "Anvil"
$deck.Range("B18").Value = 1
$deck.Range("C18").Value = 3
$deck.Range("D18").Value = "Steel"
$deck.Range("H18").Value = "anvil"

# Row 19: Baseball Bat
$deck.Range("A19").Value = "Baseball Bat"
$deck.Range("B19").Value = 1
$deck.Range("C19").Value = 3
$deck.Range("D19").Value = "Wood"
$deck.Range("H19").Value = "baseball-bat"

# Row 20: Crossbow
$deck.Range("A20").Value = "Crossbow"
$deck.Range("B20").Value = 1
$deck.Range("C20").Value = 2
$deck.Range("D20").Value = "Wood"
$deck.Range("E20").Value = 2
$deck.Range("F20").Value = "Steel"
$deck.Range("H20").Value = "crossbow"

# --- 5. VP formula for every card row (6-27) -----------------------------
for ($r = 6; $r -le 27; $r++) {
    $formula = "=ROUND(C$r*VLOOKUP(D$r,VPs!A`$2:B`$5,2,FALSE) + IF(ISNA(VLOOKUP(F$r,VPs!A`$2:B`$6,2,FALSE)),0,E$r*VLOOKUP(F$r,VPs!A`$2:B`$6,2,FALSE)),0)"
    $deck.Range("G$r").Formula = $formula
}

# --- 6. Selections to match the final saved state ------------------------
$deck.Activate()
$deck.Range("C26").Select()
$vps.Range("B6").Select()
$deck.Activate()
